# Applies the cryptos.xlsx data refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) hold text-like values (thousand-dot
# separators, padded percentages) that Excel would otherwise reinterpret
# as numbers. Mark the range as Text first, write the values, then restore
# the default style so the saved cells match the original (unstyled) cells.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '66.603.12'
$ws.Range('E2').Value = '  -4.42%  '
$ws.Range('D3').Value = '3.368.90'
$ws.Range('E3').Value = '  -5.11%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = '561.16'
$ws.Range('E5').Value = '  -4.45%  '
$ws.Range('D6').Value = '183.25'
$ws.Range('E6').Value = '  -7.70%  '
$ws.Range('D7').Value = '0.600'
$ws.Range('E7').Value = '  -2.43%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').Value = '3.359.98'
$ws.Range('E9').Value = '  -5.01%  '
$ws.Range('E10').Value = '  -9.94%  '
$ws.Range('E11').Value = '  -5.51%  '
$ws.Range('D12').Value = '48.15'
$ws.Range('E12').Value = '  -8.04%  '
$ws.Range('E13').Value = '  -7.40%  '
$ws.Range('E14').Value = '  -6.40%  '
$ws.Range('D15').Value = '3.914.29'
$ws.Range('E15').Value = '  -4.74%  '
$ws.Range('D16').Value = '607.39'
$ws.Range('E16').Value = '  -11.95%  '
$ws.Range('D17').Value = '18.30'
$ws.Range('E17').Value = '  -1.97%  '
$ws.Range('D18').Value = '66.488.88'
$ws.Range('E18').Value = '  -4.62%  '
$ws.Range('D19').Value = '3.374.16'
$ws.Range('E19').Value = '  -5.29%  '
$ws.Range('E20').Value = '  -3.22%  '
$ws.Range('D21').Value = '11.58'
$ws.Range('E21').Value = '  -7.45%  '
$ws.Range('D22').Value = '0.917'
$ws.Range('E22').Value = '  -5.95%  '
$ws.Range('D23').Value = '17.03'
$ws.Range('E23').Value = '  -7.49%  '
$ws.Range('E24').Value = '  -1.06%  '
$ws.Range('D25').Value = '98.81'
$ws.Range('E25').Value = '  -9.16%  '
$ws.Range('D26').Value = '4.08'
$ws.Range('E26').Value = '  -7.73%  '
$ws.Range('D27').Value = '6.00'
$ws.Range('E28').Value = '  -8.28%  '
$ws.Range('D29').Value = '9.44'
$ws.Range('E29').Value = '  -8.39%  '
$ws.Range('E30').Value = '  -10.02%  '
$ws.Range('D31').Value = '30.81'
$ws.Range('E31').Value = '  -8.69%  '
$ws.Range('D32').Value = '6.33'
$ws.Range('E32').Value = '  -8.69%  '
$ws.Range('D33').Value = '3.84'
$ws.Range('E33').Value = '  -12.96%  '
$ws.Range('E34').Value = '  -6.54%  '
$ws.Range('D35').Value = '552.53'
$ws.Range('E35').Value = '  +10.69%  '
$ws.Range('E36').Value = '  -5.40%  '
$ws.Range('D37').Value = '3.813.09'
$ws.Range('E37').Value = '  +0.12%  '
$ws.Range('D38').Value = '58.24'
$ws.Range('E38').Value = '  -6.33%  '
$ws.Range('D39').Value = '0.998'
$ws.Range('E39').Value = '  -0.10%  '
$ws.Range('D40').Value = '3.43'
$ws.Range('E40').Value = '  -7.30%  '
$ws.Range('E41').Value = '  -11.97%  '
$ws.Range('B42').Value = 'Fetch.AI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D42').Value = '2.70'
$ws.Range('E42').Value = '  -8.86%  '
$ws.Range('B43').Value = 'CoreDAO'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D43').Value = '3.41'
$ws.Range('E43').Value = '  +21.48%  '
$ws.Range('E44').Value = '  -5.87%  '
$ws.Range('E45').Value = '  -6.74%  '
$ws.Range('D46').Value = '32.50'
$ws.Range('E46').Value = '  -7.00%  '
$ws.Range('D47').Value = '0.0421'
$ws.Range('E47').Value = '  -9.24%  '
$ws.Range('D48').Value = '3.22'
$ws.Range('E48').Value = '  -4.84%  '
$ws.Range('D49').Value = '2.67'
$ws.Range('E49').Value = '  -9.34%  '
$ws.Range('E50').Value = '  -5.18%  '
$ws.Range('E51').Value = '  -0.06%  '

$ws.Range("D2:E51").Style = "Normal"
